$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "47×81="
$t.Cell(1, 2).Range.Text = "64×95="
$t.Cell(1, 3).Range.Text = "69×44="
$t.Cell(1, 4).Range.Text = "24×70="
$t.Cell(1, 5).Range.Text = "22×29="

$t.Cell(5, 1).Range.Text = "23×47="
$t.Cell(5, 2).Range.Text = "13×93="
$t.Cell(5, 3).Range.Text = "70×35="
$t.Cell(5, 4).Range.Text = "70×41="
$t.Cell(5, 5).Range.Text = "29×58="

$t.Cell(10, 1).Range.Text = "49×29="
$t.Cell(10, 2).Range.Text = "35×93="
$t.Cell(10, 3).Range.Text = "31×97="
$t.Cell(10, 4).Range.Text = "49×20="
$t.Cell(10, 5).Range.Text = "85×41="

$t.Cell(15, 1).Range.Text = "75×29="
$t.Cell(15, 2).Range.Text = "82×79="
$t.Cell(15, 3).Range.Text = "41×53="
$t.Cell(15, 4).Range.Text = "16×56="
$t.Cell(15, 5).Range.Text = "11×64="

$t.Cell(20, 1).Range.Text = "56×28="
$t.Cell(20, 2).Range.Text = "37×33="
$t.Cell(20, 3).Range.Text = "71×41="
$t.Cell(20, 4).Range.Text = "61×57="
$t.Cell(20, 5).Range.Text = "25×48="
